$d = $word.ActiveDocument

# 1) Update the document title text.
$d.Content.Find.Execute(
    "Lista de tecnologias utilizadas na solução", $true, $false, $false,
    $false, $false, $true, 1, $false, "Lista de Tecnologias", 2) | Out-Null

# 2) Insert two new bullet paragraphs ("- Heroku Application Service: ..."
#    and "- GitHub: ...") right after the "- Maven: ..." paragraph, before
#    the blank paragraph that precedes the "2. Database" heading.

# Locate the "- Maven: ..." paragraph; the blank paragraph right after it
# is where the new content needs to land (inserted *before* that blank
# paragraph so it keeps clean/inherited-free formatting).
$mavenIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -match "^- Maven:") {
        $mavenIndex = $i
        break
    }
}

$blank = $d.Paragraphs.Item($mavenIndex + 1)

# --- "- Heroku Application Service: " paragraph ---
$blank.Range.InsertParagraphBefore()
$heroku = $d.Paragraphs.Item($mavenIndex + 1)

$boldPart = "- Heroku Application Service: "
$restPart = "Utilizado como servidor de aplicação cloud, onde hospedaremos nossa aplicação Java."
$heroku.Range.Text = $boldPart + $restPart

$start = $heroku.Range.Start
$d.Range($start, $start + $boldPart.Length).Font.Bold = 1

# --- "- GitHub: " paragraph ---
$blank2 = $d.Paragraphs.Item($mavenIndex + 2)
$blank2.Range.InsertParagraphBefore()
$github = $d.Paragraphs.Item($mavenIndex + 2)

$boldPart2 = "- GitHub: "
$restPart2 = "Utilizado como versionador de software, onde mantemos o controle do código fonte e suas versões durante o desenvolvimento do projeto, assim como o controle dos responsáveis pelo desenvolvimento."
$github.Range.Text = $boldPart2 + $restPart2

$start2 = $github.Range.Start
$d.Range($start2, $start2 + $boldPart2.Length).Font.Bold = 1
